# Fix overlap token during NER training, add poetry and update documentation
#
# - "train" sheet: the "Uber blew through $1 million a week" training
#   example (plus its two entity rows "Uber"/ORG and "$1 million"/MONEY)
#   overlapped another entity and is removed.
# - Header rows across the sheets are made bold.
# - "train_iteration" on the config sheet goes from 2 to 20.
# - Selection / active-sheet bookkeeping moves to reflect where the
#   author ended up working (source!C2, train!B13, config!B5, with the
#   config sheet becoming the active tab).

$wb = $excel.ActiveWorkbook

$wsSource  = $wb.Worksheets.Item("source")
$wsPrepare = $wb.Worksheets.Item("prepare")
$wsTrain   = $wb.Worksheets.Item("train")
$wsConfig  = $wb.Worksheets.Item("config")

# --- train: remove the overlapping "Uber blew through $1 million a week"
#     example and its two entity rows (old rows 2-4) ---------------------
$wsTrain.Range("A2:A4").EntireRow.Delete()

# --- bold the header row on every sheet ---------------------------------
$wsSource.Range("A1:C1").Font.Bold = $true
$wsPrepare.Range("A1:D1").Font.Bold = $true
$wsTrain.Range("A1:D1").Font.Bold = $true
$wsConfig.Range("A1:B1").Font.Bold = $true

# --- config: bump train_iteration from 2 to 20 --------------------------
$wsConfig.Range("B4").Value = 20

# --- selections on each sheet, ending with config active ----------------
[void]$wsSource.Activate()
[void]$wsSource.Range("C2").Select()

[void]$wsTrain.Activate()
[void]$wsTrain.Range("B13").Select()

[void]$wsConfig.Activate()
[void]$wsConfig.Range("B5").Select()

Write-Output "edits applied"
